# Generate Report for Handoff
# - Flip the handback status back to "Ready for handoff" (report was
#   regenerated for a fresh handoff instead of reporting the prior handback).
# - Refresh the "generate date" / "handoff datetime" timestamps that go with
#   that status on each language sheet (and the rollup on Overview).
# - Status column got noticeably shorter text, so the report's column
#   autosizing made columns E/F (Overview) and C (zh-cn, de-de) narrower.

$wb = $excel.ActiveWorkbook

$statusOld = "Handed back: in sync with en-US"
$statusNew = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status + the Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("G2").Value = "2016-11-08 23:33:16"

$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# ---------------------------------------------------------------------
# zh-cn sheet: Status + Latest Handoff Datetime
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("H2").Value = "2016-11-08 23:33:03"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# ---------------------------------------------------------------------
# de-de sheet: Status + Latest Handoff Datetime
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("H2").Value = "2016-11-08 23:33:16"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
